$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19 held leftover test data ("Teste Hugo" / "Alberto" ...). Remove it.
# Columns whose explicit cell style equals the column's own default style
# are fully cleared (value + formatting), while the remaining columns keep
# their (non-default) formatting and only lose their value/hyperlink.
$ws.Range("A19:D19").Clear()
$ws.Range("F19:F19").Clear()
$ws.Range("H19:R19").Clear()
$ws.Range("X19:Y19").Clear()
$ws.Range("AE19:AF19").Clear()
$ws.Range("AL19:AM19").Clear()
$ws.Range("AZ19:BA19").Clear()

$ws.Range("E19:E19").ClearContents()
$ws.Range("G19:G19").ClearContents()
$ws.Range("S19:W19").ClearContents()
$ws.Range("Z19:AD19").ClearContents()
$ws.Range("AG19:AK19").ClearContents()
$ws.Range("AN19:AY19").ClearContents()
$ws.Range("BB19:BH19").ClearContents()

# Leave the selection on the cleared row, matching the saved view state.
$ws.Activate()
$ws.Range("A19:XFD19").Select()
